# Scheduled runner update: refresh market-price-derived figures on several
# per-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) of Carbuncle_Profits.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 326.5
$ws.Range("J2").Value = 382.4
$ws.Range("L2").Value = 382.4
$ws.Range("N2").Value = -608.4
$ws.Range("H15").Value = 17847766
$ws.Range("I15").Value = 17847766
$ws.Range("K15").Value = 53543298
$ws.Range("M15").Value = -53543129
$ws.Range("H43").Value = 1145
$ws.Range("I43").Value = 900
$ws.Range("J43").Value = 1226.6666
$ws.Range("K43").Value = 900
$ws.Range("L43").Value = 1226.6666
$ws.Range("M43").Value = -831
$ws.Range("N43").Value = -1364.6666
$ws.Range("H46").Value = 20833.65
$ws.Range("J46").Value = 20833.65
$ws.Range("L46").Value = 62500.95
$ws.Range("N46").Value = -62738.95
$ws.Range("H60").Value = 20833.65
$ws.Range("J60").Value = 20833.65
$ws.Range("L60").Value = 62500.95
$ws.Range("N60").Value = -63468.95
$ws.Range("H100").Value = 45643460
$ws.Range("I100").Value = 64816660
$ws.Range("J100").Value = 2503751.5
$ws.Range("K100").Value = 64816660
$ws.Range("L100").Value = 2503751.5
$ws.Range("M100").Value = -64816119
$ws.Range("N100").Value = -2504833.5
$ws.Range("H114").Value = 37806.855
$ws.Range("J114").Value = 37806.855
$ws.Range("L114").Value = 37806.855
$ws.Range("N114").Value = -46484.855
$ws.Range("H121").Value = 1135.625
$ws.Range("I121").Value = 395
$ws.Range("J121").Value = 1241.4286
$ws.Range("K121").Value = 1185
$ws.Range("L121").Value = 3724.2858
$ws.Range("M121").Value = 562
$ws.Range("N121").Value = -7218.2858
$ws.Range("H125").Value = 1308084.2
$ws.Range("I125").Value = 473.33334
$ws.Range("J125").Value = 1961889.6
$ws.Range("K125").Value = 4260.00006
$ws.Range("L125").Value = 17657006.4
$ws.Range("M125").Value = -1800.00006
$ws.Range("N125").Value = -17661926.4
$ws.Range("H132").Value = 5158
$ws.Range("I132").Value = 5824.25
$ws.Range("J132").Value = 2493
$ws.Range("K132").Value = 17472.75
$ws.Range("L132").Value = 7479
$ws.Range("M132").Value = -14942.75
$ws.Range("N132").Value = -12539
$ws.Range("H135").Value = 612.7917
$ws.Range("I135").Value = 246.75
$ws.Range("J135").Value = 2443
$ws.Range("K135").Value = 2220.75
$ws.Range("L135").Value = 21987
$ws.Range("M135").Value = 314.25
$ws.Range("N135").Value = -27057
$ws.Range("H138").Value = 3413.4595
$ws.Range("J138").Value = 5766
$ws.Range("L138").Value = 17298
$ws.Range("N138").Value = -27578
$ws.Range("H141").Value = 7337.1113
$ws.Range("I141").Value = 2607
$ws.Range("J141").Value = 13249.75
$ws.Range("K141").Value = 7821
$ws.Range("L141").Value = 39749.25
$ws.Range("M141").Value = -2641
$ws.Range("N141").Value = -50109.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9361.714
$ws.Range("I32").Value = 5152.706
$ws.Range("K32").Value = 5152.706
$ws.Range("M32").Value = -4865.706
$ws.Range("H74").Value = 2040.0303
$ws.Range("I74").Value = 2150.889
$ws.Range("J74").Value = 1907
$ws.Range("K74").Value = 2150.889
$ws.Range("L74").Value = 1907
$ws.Range("M74").Value = -1276.889
$ws.Range("N74").Value = -3655
$ws.Range("H77").Value = 2040.0303
$ws.Range("I77").Value = 2150.889
$ws.Range("J77").Value = 1907
$ws.Range("K77").Value = 10754.445
$ws.Range("L77").Value = 9535
$ws.Range("M77").Value = -6386.445
$ws.Range("N77").Value = -18271
$ws.Range("H132").Value = 3053.037
$ws.Range("I132").Value = 1285.8462
$ws.Range("J132").Value = 4694
$ws.Range("K132").Value = 3857.5386
$ws.Range("L132").Value = 14082
$ws.Range("M132").Value = -1327.5386
$ws.Range("N132").Value = -19142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 54905
$ws.Range("J55").Value = 54905
$ws.Range("L55").Value = 54905
$ws.Range("N55").Value = -55451
$ws.Range("H105").Value = 3049.625
$ws.Range("I105").Value = 3049.625
$ws.Range("K105").Value = 3049.625
$ws.Range("M105").Value = -1302.625
$ws.Range("H107").Value = 729.5333000000001
$ws.Range("I107").Value = 687.9231
$ws.Range("K107").Value = 687.9231
$ws.Range("M107").Value = 1232.0769
$ws.Range("H134").Value = 2314.6843
$ws.Range("I134").Value = 1549.3103
$ws.Range("J134").Value = 4780.8887
$ws.Range("K134").Value = 4647.9309
$ws.Range("L134").Value = 14342.6661
$ws.Range("M134").Value = -2112.9309
$ws.Range("N134").Value = -19412.6661

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14169.207
$ws.Range("I31").Value = 1517.5238
$ws.Range("J31").Value = 47379.875
$ws.Range("K31").Value = 1517.5238
$ws.Range("L31").Value = 47379.875
$ws.Range("M31").Value = -1222.5238
$ws.Range("N31").Value = -47969.875
$ws.Range("H34").Value = 14169.207
$ws.Range("I34").Value = 1517.5238
$ws.Range("J34").Value = 47379.875
$ws.Range("K34").Value = 1517.5238
$ws.Range("L34").Value = 47379.875
$ws.Range("M34").Value = -1315.5238
$ws.Range("N34").Value = -47783.875
$ws.Range("H132").Value = 1819.931
$ws.Range("I132").Value = 1456.9231
$ws.Range("J132").Value = 4966
$ws.Range("K132").Value = 4370.7693
$ws.Range("L132").Value = 14898
$ws.Range("M132").Value = -1840.7693
$ws.Range("N132").Value = -19958

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 29
$ws.Range("J2").Value = 32.8
$ws.Range("L2").Value = 196.8
$ws.Range("N2").Value = -422.8
$ws.Range("H5").Value = 2107.375
$ws.Range("I5").Value = 1126.6666
$ws.Range("J5").Value = 2333.6924
$ws.Range("K5").Value = 3379.9998
$ws.Range("L5").Value = 7001.0772
$ws.Range("M5").Value = -3267.9998
$ws.Range("N5").Value = -7225.0772
$ws.Range("H23").Value = 134.3
$ws.Range("I23").Value = 129
$ws.Range("J23").Value = 134.88889
$ws.Range("K23").Value = 387
$ws.Range("L23").Value = 404.66667
$ws.Range("M23").Value = -152
$ws.Range("N23").Value = -874.6666700000001
$ws.Range("H120").Value = 3252.4443
$ws.Range("I120").Value = 2034.1428
$ws.Range("J120").Value = 7516.5
$ws.Range("K120").Value = 6102.428400000001
$ws.Range("L120").Value = 22549.5
$ws.Range("M120").Value = -1264.428400000001
$ws.Range("N120").Value = -32225.5
$ws.Range("H132").Value = 1059.1578
$ws.Range("I132").Value = 889
$ws.Range("K132").Value = 8001
$ws.Range("M132").Value = -5471
$ws.Range("H133").Value = 2153.3333
$ws.Range("I133").Value = 2328.5715
$ws.Range("K133").Value = 6985.7145
$ws.Range("M133").Value = -1925.7145
$ws.Range("H134").Value = 1500.3214
$ws.Range("I134").Value = 1067.2667
$ws.Range("K134").Value = 3201.800099999999
$ws.Range("M134").Value = 1868.199900000001
$ws.Range("H135").Value = 2107.375
$ws.Range("I135").Value = 1126.6666
$ws.Range("J135").Value = 2333.6924
$ws.Range("K135").Value = 10139.9994
$ws.Range("L135").Value = 21003.2316
$ws.Range("M135").Value = -7604.999400000001
$ws.Range("N135").Value = -26073.2316

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3951.6667
$ws.Range("I132").Value = 2290.889
$ws.Range("J132").Value = 5612.4443
$ws.Range("K132").Value = 6872.667
$ws.Range("L132").Value = 16837.3329
$ws.Range("M132").Value = -4342.667
$ws.Range("N132").Value = -21897.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 12553
$ws.Range("J80").Value = 12553
$ws.Range("L80").Value = 12553
$ws.Range("N80").Value = -14799
$ws.Range("H83").Value = 12553
$ws.Range("J83").Value = 12553
$ws.Range("L83").Value = 37659
$ws.Range("N83").Value = -48891
$ws.Range("H132").Value = 24395200
$ws.Range("I132").Value = 40005024
$ws.Range("J132").Value = 4848.3125
$ws.Range("K132").Value = 120015072
$ws.Range("L132").Value = 14544.9375
$ws.Range("M132").Value = -120012542
$ws.Range("N132").Value = -19604.9375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4414.0815
$ws.Range("I136").Value = 796.6667
$ws.Range("K136").Value = 2390.0001
$ws.Range("M136").Value = 159.9998999999998
